# "Add files via upload" — replace the trailing blank placeholder rows with
# four new wishlist items, rename the sheet, and shrink the autofilter /
# sort-state ranges down to the new data extent (row order is NOT changed —
# the new rows are simply appended after the existing, unsorted, ones).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Extra Wishlist"

# Clear out the old empty placeholder rows (16:23) before writing new data.
$ws.Range("A16:E23").Clear()

# New wishlist entries.
$ws.Range("A16").Value = "Sigma 24-70mm f/2.8 DG DN II Art"
$ws.Range("B16").Value = "https://static01.galaxus.com/productimages/3/4/1/4/6/9/1/6/1/2/1/6/6/2/2/1/8/4/9/682182dc-69c0-4397-b129-6fef5b7c1b49_cropped.jpg_2880.avif"
$ws.Range("C16").Value = "https://www.digitec.ch/en/s1/product/sigma-24-70mm-f28-dg-dn-ii-art-sony-e-full-size-lenses-45881601"
$ws.Range("D16").Value = 1160

$ws.Range("A17").Value = "Deuter Forest Fox 10"
$ws.Range("B17").Value = "https://static01.galaxus.com/productimages/3/3/3/3/4/0/2/0/0/1/7/4/8/9/2/2/0/5/2/58cd5acc-8462-46b1-8934-2a09baf72ca7_cropped.jpg_720.avif"
$ws.Range("C17").Value = "https://www.galaxus.ch/en/s8/product/deuter-forest-fox-10-10-l-backpacks-42902192"
$ws.Range("C17").Style = "Hyperlink"
$ws.Range("D17").Value = 45

$ws.Range("A18").Value = "Deuter Forest Fox 10"
$ws.Range("B18").Value = "https://static01.galaxus.com/productimages/1/6/1/0/9/8/0/5/7/2/6/1/7/7/2/6/8/4/5/118d9d2b-617e-4c6c-ba5b-b3e0fac55598_cropped.jpg_720.avif"
$ws.Range("C18").Value = "https://www.galaxus.ch/en/s8/product/deuter-waldfuchs-10-10-l-backpacks-42902191"
$ws.Range("C18").Style = "Hyperlink"
$ws.Range("D18").Value = 45

$ws.Range("A19").Value = "Sigma 100-400mm f/5.0-6.3 DG DN OS, Sony E"
$ws.Range("B19").Value = "https://static01.galaxus.com/productimages/3/5/9/9/6/4/4/5/PPhoto_100_400_5_6.3_dgdn_c020_Lmt_horizontal_tripodsocket.jpg_720.avif"
$ws.Range("C19").Value = "https://www.digitec.ch/en/s1/product/sigma-100-400mm-f50-63-dg-dn-os-sony-e-sony-e-full-size-lenses-13366507"
$ws.Range("C19").Style = "Hyperlink"
$ws.Range("D19").Value = 870

# Selection moves as Excel would leave it after entering the data.
$ws.Range("E26").Select()

# Shrink the autofilter range to the new data extent.
$ws.AutoFilterMode = $false
$ws.Range("A1:E19").AutoFilter()

# Shrink the sort-state ranges to match. Applying a real sort always
# physically reorders rows, but the target data stays in plain insertion
# order (it's not actually sorted by price). So: make column D temporarily
# strictly increasing (already-sorted w.r.t. an ascending sort => the Apply
# below is a no-op on row order), apply the ascending sort to get Excel to
# persist the shrunk sortState/sortCondition refs, then write the real
# price values back over the top.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value = $r
}
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D2:D19"))
$ws.Sort.SetRange($ws.Range("A2:D19"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

$prices = @(1600, 750, 1200, 3800, 700, 540, 250, 750, 330, 220, 240, 160, 240, 1700, 1160, 45, 45, 870)
for ($i = 0; $i -lt $prices.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $prices[$i]
}

# AutoFilter doesn't retarget the hidden _FilterDatabase defined name on its
# own — update it explicitly to match the new range/sheet name.
$wb.Names.Item(1).RefersTo = "='Extra Wishlist'!`$A`$1:`$E`$19"

$wb.Save()
